$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.91
$ws.Range("I2").Value = 4.1
$ws.Range("U2").Value = 9
$ws.Range("X2").Value = 15
$ws.Range("AD2").Value = 251
$ws.Range("AG2").Value = 15
$ws.Range("G3").Value = 1.98
$ws.Range("H3").Value = 2.92
$ws.Range("L3").Value = 1.4
$ws.Range("M3").Value = 2.52
$ws.Range("N3").Value = 2.18
$ws.Range("O3").Value = 1.53
$ws.Range("P3").Value = 1.45
$ws.Range("Q3").Value = 2.37
$ws.Range("R3").Value = 1.9
$ws.Range("S3").Value = 1.72
$ws.Range("T3").Value = 5.7
$ws.Range("U3").Value = 8.5
$ws.Range("V3").Value = 8.5
$ws.Range("W3").Value = 18
$ws.Range("Y3").Value = 35
$ws.Range("AA3").Value = 5.8
$ws.Range("AB3").Value = 15.5
$ws.Range("AC3").Value = 90
$ws.Range("AD3").Value = 800
$ws.Range("AE3").Value = 10.25
$ws.Range("AF3").Value = 23
$ws.Range("AG3").Value = 13.5
$ws.Range("AH3").Value = 75
$ws.Range("AJ3").Value = 50
$ws.Range("G6").Value = 2.22
$ws.Range("H6").Value = 3.5
$ws.Range("I6").Value = 2.67
$ws.Range("T6").Value = 10
$ws.Range("U6").Value = 11.75
$ws.Range("W6").Value = 19.5
$ws.Range("Z6").Value = 15.5
$ws.Range("AA6").Value = 6.5
$ws.Range("AB6").Value = 9.5
$ws.Range("AC6").Value = 27
$ws.Range("AE6").Value = 11
$ws.Range("AF6").Value = 14.5
$ws.Range("AG6").Value = 8.75
$ws.Range("AH6").Value = 26
$ws.Range("AI6").Value = 16
$ws.Range("G7").Value = 1.45
$ws.Range("H7").Value = 4.3
$ws.Range("I7").Value = 5.4
$ws.Range("T7").Value = 8
$ws.Range("U7").Value = 7.3
$ws.Range("W7").Value = 9.25
$ws.Range("Y7").Value = 16
$ws.Range("Z7").Value = 16
$ws.Range("AB7").Value = 12.5
$ws.Range("AC7").Value = 40
$ws.Range("AE7").Value = 16
$ws.Range("AG7").Value = 14.5
$ws.Range("AI7").Value = 37
$ws.Range("AJ7").Value = 32
$ws.Range("N10").Value = 1.95
$ws.Range("O10").Value = 1.85
$ws.Range("P10").Value = 1.36
$ws.Range("Q10").Value = 3
$ws.Range("G11").Value = 2.55
$ws.Range("H11").Value = 3.5
$ws.Range("L11").Value = 1.25
$ws.Range("M11").Value = 3.75
$ws.Range("N11").Value = 1.88
$ws.Range("O11").Value = 1.93
$ws.Range("P11").Value = 1.36
$ws.Range("Q11").Value = 3
$ws.Range("R11").Value = 1.67
$ws.Range("S11").Value = 2.1
$ws.Range("T11").Value = 9
$ws.Range("Z11").Value = 11
$ws.Range("AB11").Value = 13
$ws.Range("AC11").Value = 41
$ws.Range("AD11").Value = 201
$ws.Range("AE11").Value = 9.5
$ws.Range("G13").Value = 2.15
$ws.Range("H13").Value = 3.5
$ws.Range("I13").Value = 3.25
$ws.Range("J13").Value = 1.03
$ws.Range("K13").Value = 13
$ws.Range("L13").Value = 1.22
$ws.Range("N13").Value = 1.74
$ws.Range("O13").Value = 1.94
$ws.Range("W13").Value = 19
$ws.Range("Y13").Value = 23
$ws.Range("Z13").Value = 12
$ws.Range("AD13").Value = 151
$ws.Range("AJ13").Value = 29
$ws.Range("N14").Value = 1.7
$ws.Range("O14").Value = 2.1
